$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# translation rows to append: key, value
$rows = @(
    @("error.Duplicate entry [z_build_name_unique] of [z_build].", "Tento název buildu jste již použili; zvolte prosím jiný."),
    @("lab.build.table.name", "Jméno buildu"),
    @("lab.build.table.atomizer", "Atomizér"),
    @("lab.build.table.cotton", "Vata"),
    @("lab.build.table.coil", "Spirálka"),
    @("lab.build.table.ohm", "Odpor buildu"),
    @("lab.build.table.coils", "Počet spirálek"),
    @("lab.build.table.created", "Vytvořen"),
    @("lab.build.created.message", "Build [{{data.name}}] byl uložen."),
    @("lab.build.table.coilOffset", "Pozice spirálky"),
    @("lab.build.table.cottonOffset", "Množství vaty")
)

$lastRow = 286
$startRow = $lastRow + 1
$r = $startRow

foreach ($row in $rows) {
    # Copy formatting from the last existing data row so new rows keep the
    # same cell style as the rest of the table.
    $ws.Range("A$lastRow`:C$lastRow").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]

    $r = $r + 1
}

$excel.CutCopyMode = 0

$lastNewRow = $r - 1
$ws.Range("B$lastNewRow").Select()
